# Update cryptocurrency price/volume figures per latest refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.032.38"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "1.721.50"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.65"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.07"
$ws.Range("E8").Value = "  +12.88%  "
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "1.968.39"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "1.719.66"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.26"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.559"
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.41"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "28.000.62"
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.19"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "0.0₃0754"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.62"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.65"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.60"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.48"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.68"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0509"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "1.491.53"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.26"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.951"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.604"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.47"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.81"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.30"
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").Value = "1.870.93"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.795"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.77"
$ws.Range("E47").Value = "  +12.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "90.76"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "0.0₆0112"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("E51").Value = "  +0.05%  "
